# Generate Report for Handoff
# The source file was re-handed-off under a new GUID
# (0b4c3046-7a2c-4471-a243-7a24bb35215f -> 1b594f1f-92d9-4483-aae0-e2a9ce1be894),
# with a freshly generated xliff hash (e37391fc87edcace882f353b43e03e8cfdd7e087 ->
# 9f3e8b70584138be85c8ea6331652932cf51b292). This refreshes the localization
# status report: new handoff file names/timestamps, and the not-yet-handed-back
# target/handback columns are cleared out for both locales.

$wb = $excel.ActiveWorkbook

$oldGuid = '0b4c3046-7a2c-4471-a243-7a24bb35215f'
$newGuid = '1b594f1f-92d9-4483-aae0-e2a9ce1be894'
$oldHash = 'e37391fc87edcace882f353b43e03e8cfdd7e087'
$newHash = '9f3e8b70584138be85c8ea6331652932cf51b292'

$newFileName = "$newGuid.md"
$newRelPath  = "e2e\$newGuid.md"

$newXliffZh = "$newGuid.$newHash.zh-cn.xlf"
$newXliffDe = "$newGuid.$newHash.de-de.xlf"

$neverHandedBack = '0001-01-01 00:00:00'

$githubBaseUrl = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84bc305b75ba80199bf4a9ac6c369540a47ce1df/'

# ---------------------------------------------------------------------------
# Overview sheet: new source-file name/path + refreshed "Latest HO Xliff
# Generate Date".
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newRelPath
$wsOverview.Range("G2").Value = '2016-08-17 06:53:22'

$rB2 = $wsOverview.Range("B2")
$rB2.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($rB2, $githubBaseUrl + $newRelPath, "", "", $newRelPath)

# ---------------------------------------------------------------------------
# zh-cn sheet: new source-file name, new handoff xliff + datetime, and the
# target/handback columns reset because the new handoff hasn't come back yet.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = $newXliffZh
$wsZh.Range("H2").Value = '2016-08-17 06:53:17'
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $neverHandedBack

$rZhA2 = $wsZh.Range("A2")
$rZhA2.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($rZhA2, $githubBaseUrl + $newRelPath, "", "", $newFileName)

$rZhI2 = $wsZh.Range("I2")
$rZhI2.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# de-de sheet: same story as zh-cn.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = $newXliffDe
$wsDe.Range("H2").Value = '2016-08-17 06:53:22'
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $neverHandedBack

$rDeA2 = $wsDe.Range("A2")
$rDeA2.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($rDeA2, $githubBaseUrl + $newRelPath, "", "", $newFileName)

$rDeI2 = $wsDe.Range("I2")
$rDeI2.Hyperlinks.Delete()
